# Applies the "gh-pages output generated at 456a3b4" data refresh:
# numeric bumps to the "想去人数" (interest count) column F (and a couple of
# "最低票价" column G conversions from an inlineStr "已售罄" marker to an
# actual lowest-price number), plus two corrected time-range strings in
# column E, across the 展览 / 演出 / 本地生活 / 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 37678
$ws1.Range("G2").Value = 68
$ws1.Range("F5").Value = 774
$ws1.Range("F6").Value = 483
$ws1.Range("F7").Value = 369
$ws1.Range("F8").Value = 468
$ws1.Range("F11").Value = 718
$ws1.Range("F12").Value = 556
$ws1.Range("F13").Value = 50
$ws1.Range("F14").Value = 36
$ws1.Range("F15").Value = 24
$ws1.Range("F16").Value = 653
$ws1.Range("F17").Value = 180
$ws1.Range("F18").Value = 472
$ws1.Range("F19").Value = 443
$ws1.Range("F20").Value = 1170
$ws1.Range("E21").Value = "2024.05.10 10:00-05.12 18:00"
$ws1.Range("F21").Value = 93
$ws1.Range("F22").Value = 835
$ws1.Range("F23").Value = 2539
$ws1.Range("F24").Value = 1016
$ws1.Range("F25").Value = 563
$ws1.Range("F26").Value = 109
$ws1.Range("F29").Value = 780
$ws1.Range("F30").Value = 62
$ws1.Range("F31").Value = 1162

# --- Sheet 2: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 406
$ws2.Range("F4").Value = 332

# --- Sheet 3: 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 637

# --- Sheet 4: 全部类型 (All types, union of the above) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 637
$ws4.Range("F3").Value = 37679
$ws4.Range("G3").Value = 68
$ws4.Range("F6").Value = 774
$ws4.Range("F7").Value = 483
$ws4.Range("F9").Value = 369
$ws4.Range("F10").Value = 468
$ws4.Range("F11").Value = 406
$ws4.Range("F12").Value = 332
$ws4.Range("F17").Value = 718
$ws4.Range("F18").Value = 556
$ws4.Range("F19").Value = 50
$ws4.Range("F21").Value = 36
$ws4.Range("F25").Value = 24
$ws4.Range("F27").Value = 653
$ws4.Range("F28").Value = 180
$ws4.Range("F29").Value = 472
$ws4.Range("F30").Value = 443
$ws4.Range("F31").Value = 1170
$ws4.Range("E32").Value = "2024.05.10 10:00-05.12 18:00"
$ws4.Range("F32").Value = 93
$ws4.Range("F33").Value = 835
$ws4.Range("F34").Value = 2539
$ws4.Range("F35").Value = 1016
$ws4.Range("F36").Value = 563
$ws4.Range("F37").Value = 109
$ws4.Range("F41").Value = 780
$ws4.Range("F42").Value = 62
$ws4.Range("F43").Value = 1162
